$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:D1) -------------------------------------------------
$ws.Range("A1").Value = "баркод"
$ws.Range("B1").Value = "артикул"
$ws.Range("C1").Value = "СС с НДС"
$ws.Range("D1").Value = "СС без НДС"

# New header cells (C1/D1) need the same bold/border/centered look as
# A1/B1 - copy the direct formatting over from A1.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column A (barcodes) must stay text, not be coerced to numbers -----
$a2 = $ws.Range("A2")
$a2.NumberFormat = "@"
$a2.Value = "4600956001999"
$a2.Style = "Normal"

$a3 = $ws.Range("A3")
$a3.NumberFormat = "@"
$a3.Value = "4600956006390"
$a3.Style = "Normal"

$a4 = $ws.Range("A4")
$a4.NumberFormat = "@"
$a4.Value = "4600956009759"
$a4.Style = "Normal"

$a5 = $ws.Range("A5")
$a5.NumberFormat = "@"
$a5.Value = "4600956001777"
$a5.Style = "Normal"

$a6 = $ws.Range("A6")
$a6.NumberFormat = "@"
$a6.Value = "4600956002040"
$a6.Style = "Normal"

# --- Column B (article / product names) ---------------------------------
$ws.Range("B2").Value = "Хлопья_арахис 250г"
$ws.Range("B3").Value = "Ассорти3,0"
$ws.Range("B4").Value = "Ирис_молочн_байтс500г"
$ws.Range("B5").Value = "Хлопья_мед 300г"
$ws.Range("B6").Value = "Конфеты_микс_лимт 100г"

# --- Column C (cost price with VAT) --------------------------------------
$ws.Range("C2").Value = 92.983
$ws.Range("C3").Value = 1071.18
$ws.Range("C4").Value = 186.696
$ws.Range("C5").Value = 83.688
$ws.Range("C6").Value = 187.668

# --- Column D (cost price without VAT) ------------------------------------
$ws.Range("D2").Value = 84.53
$ws.Range("D3").Value = 892.65
$ws.Range("D4").Value = 155.58
$ws.Range("D5").Value = 76.08
$ws.Range("D6").Value = 156.34

Write-Output "done"
